# Update the "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output values.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    4  = 6233
    5  = 171
    7  = 37
    8  = 1874
    9  = 1408
    11 = 953
    12 = 223
    13 = 5576
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
